# EPBDS-12588: trim the accidental duplicate 3rd ("D") column from the
# Datatype Package definition table (rows 13-16) in the 2-dim-array test
# fixture. The table only needs two columns (infoField / mainField); the
# third column was a leftover duplicate of column C and is removed here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13:D16").ClearContents()

# Restore the viewport selection recorded the last time the sheet was
# saved interactively.
$ws.Range("H19").Select() | Out-Null
